$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 137
$ws.Range("H137").Value = 27354426
$ws.Range("I137").Value = 9616406
$ws.Range("J137").Value = 46570612
$ws.Range("K137").Value = 28849218
$ws.Range("L137").Value = 139711836
$ws.Range("M137").Value = -28846668
$ws.Range("N137").Value = -139716936
# Row 138
$ws.Range("H138").Value = 2013.8539
$ws.Range("I138").Value = 1323.1608
$ws.Range("J138").Value = 3185.9395
$ws.Range("K138").Value = 3969.4824
$ws.Range("L138").Value = 9557.818499999999
$ws.Range("M138").Value = 1170.5176
$ws.Range("N138").Value = -19837.8185

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 556602
$ws.Range("I45").Value = 625976.3
$ws.Range("J45").Value = 1607.5
$ws.Range("K45").Value = 625976.3
$ws.Range("L45").Value = 1607.5
$ws.Range("M45").Value = -625599.3
$ws.Range("N45").Value = -2361.5
# Row 74
$ws.Range("H74").Value = 46154816
$ws.Range("I74").Value = 47619570
$ws.Range("J74").Value = 40002850
$ws.Range("K74").Value = 47619570
$ws.Range("L74").Value = 40002850
$ws.Range("M74").Value = -47618696
$ws.Range("N74").Value = -40004598
# Row 77
$ws.Range("H77").Value = 46154816
$ws.Range("I77").Value = 47619570
$ws.Range("J77").Value = 40002850
$ws.Range("K77").Value = 238097850
$ws.Range("L77").Value = 200014250
$ws.Range("M77").Value = -238093482
$ws.Range("N77").Value = -200022986

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 16
$ws.Range("H16").Value = 3000
$ws.Range("I16").Value = 3000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2830
# Row 23
$ws.Range("H23").Value = 2012
$ws.Range("I23").Value = 2012
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 2012
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -1729
# Row 107
$ws.Range("H107").Value = 794
$ws.Range("I107").Value = 794
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 794
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1126
$ws.Range("N107").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 980
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 980
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 980
$ws.Range("N2").Value = -1206
$ws.Range("M2").ClearContents()
# Row 6
$ws.Range("H6").Value = 1669633.5
$ws.Range("I6").Value = 5001000.5
$ws.Range("J6").Value = 3950
$ws.Range("K6").Value = 5001000.5
$ws.Range("L6").Value = 3950
$ws.Range("M6").Value = -5000887.5
$ws.Range("N6").Value = -4176
# Row 11
$ws.Range("H11").Value = 903
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 903
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 903
$ws.Range("N11").Value = -1183
$ws.Range("M11").ClearContents()
# Row 31
$ws.Range("H31").Value = 2978287.8
$ws.Range("I31").Value = 2127.625
$ws.Range("J31").Value = 12502000
$ws.Range("K31").Value = 2127.625
$ws.Range("L31").Value = 12502000
$ws.Range("M31").Value = -1832.625
$ws.Range("N31").Value = -12502590
# Row 33
$ws.Range("H33").Value = 405.5
$ws.Range("I33").Value = 405.5
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 405.5
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -26.5
# Row 34
$ws.Range("H34").Value = 2978287.8
$ws.Range("I34").Value = 2127.625
$ws.Range("J34").Value = 12502000
$ws.Range("K34").Value = 2127.625
$ws.Range("L34").Value = 12502000
$ws.Range("M34").Value = -1925.625
$ws.Range("N34").Value = -12502404
# Row 93
$ws.Range("H93").Value = 5604.6665
$ws.Range("I93").Value = 5604.6665
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 5604.6665
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -3732.6665
$ws.Range("N93").ClearContents()
# Row 122
$ws.Range("H122").Value = 7091.6924
$ws.Range("I122").Value = 10030.889
$ws.Range("J122").Value = 478.5
$ws.Range("K122").Value = 30092.667
$ws.Range("L122").Value = 1435.5
$ws.Range("M122").Value = -27642.667
$ws.Range("N122").Value = -6335.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 883.6405999999999
$ws.Range("I113").Value = 488.14285
$ws.Range("J113").Value = 1076.7906
$ws.Range("K113").Value = 1464.42855
$ws.Range("L113").Value = 3230.3718
$ws.Range("M113").Value = 705.5714499999999
$ws.Range("N113").Value = -7570.3718

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 49
$ws.Range("H49").Value = 12800
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 12800
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 12800
$ws.Range("N49").Value = -13168
# Row 93
$ws.Range("H93").Value = 12000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 12000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 12000
$ws.Range("N93").Value = -15744
# Row 102
$ws.Range("H102").Value = 4000
$ws.Range("I102").Value = 4604.8096
$ws.Range("J102").Value = 1459.8
$ws.Range("K102").Value = 4604.8096
$ws.Range("L102").Value = 1459.8
$ws.Range("M102").Value = -2982.8096
$ws.Range("N102").Value = -4703.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 32
$ws.Range("H32").Value = 2000
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1683
# Row 61
$ws.Range("H61").Value = 3900
$ws.Range("I61").Value = 3500
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 3500
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -3298
$ws.Range("N61").Value = -4404
# Row 82
$ws.Range("H82").Value = 4644.278
$ws.Range("I82").Value = 1690.7273
$ws.Range("J82").Value = 9285.571
$ws.Range("K82").Value = 1690.7273
$ws.Range("L82").Value = 9285.571
$ws.Range("M82").Value = -1329.7273
$ws.Range("N82").Value = -10007.571
# Row 85
$ws.Range("H85").Value = 4644.278
$ws.Range("I85").Value = 1690.7273
$ws.Range("J85").Value = 9285.571
$ws.Range("K85").Value = 1690.7273
$ws.Range("L85").Value = 9285.571
$ws.Range("M85").Value = -442.7273
$ws.Range("N85").Value = -11781.571
# Row 100
$ws.Range("H100").Value = 1622.7037
$ws.Range("I100").Value = 1239.65
$ws.Range("J100").Value = 2717.1428
$ws.Range("K100").Value = 1239.65
$ws.Range("L100").Value = 2717.1428
$ws.Range("M100").Value = -698.6500000000001
$ws.Range("N100").Value = -3799.1428
# Row 113
$ws.Range("H113").Value = 3900
$ws.Range("I113").Value = 3500
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 3500
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -1330
$ws.Range("N113").Value = -8340
# Row 132
$ws.Range("H132").Value = 4766710
$ws.Range("I132").Value = 5295604
$ws.Range("J132").Value = 6666
$ws.Range("K132").Value = 15886812
$ws.Range("L132").Value = 19998
$ws.Range("M132").Value = -15884282
$ws.Range("N132").Value = -25058

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 74
$ws.Range("H74").Value = 25750
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 25750
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 25750
$ws.Range("N74").Value = -27622
# Row 77
$ws.Range("H77").Value = 25750
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 25750
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 77250
$ws.Range("N77").Value = -86610
# Row 122
$ws.Range("H122").Value = 1440.25
$ws.Range("I122").Value = 1317
$ws.Range("J122").Value = 1810
$ws.Range("K122").Value = 3951
$ws.Range("L122").Value = 5430
$ws.Range("M122").Value = -1501
$ws.Range("N122").Value = -10330
# Row 126
$ws.Range("H126").Value = 1307.4445
$ws.Range("I126").Value = 375.30768
$ws.Range("J126").Value = 3731
$ws.Range("K126").Value = 1125.92304
$ws.Range("L126").Value = 11193
$ws.Range("M126").Value = 1344.07696
$ws.Range("N126").Value = -16133
